$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "International Financial Statis"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# 1) The second row (sub-header: "International Reserves..." / "Prices, CPI...")
#    grows from 16.5pt to 17.25pt.
$ws1.Rows.Item(2).RowHeight = 17.25

# 2) Append twelve new monthly data rows for "United States", Jan 1990..Dec 1990,
#    continuing directly after the existing last row (Dec 1989 = row 363).
$monthsData = @(
  @{ month = "Jan 1990"; v1 = 75506.361462393;  v2 = 58.4254925924822;  ht = 16.5 },
  @{ month = "Feb 1990"; v1 = 74173.5248679896; v2 = 58.7006518982553;  ht = 16.5 },
  @{ month = "Mar 1990"; v1 = 76303.0285402087; v2 = 59.0216710883238;  ht = 16.5 },
  @{ month = "Apr 1990"; v1 = 76283.2106099483; v2 = 59.1133908569149;  ht = 16.5 },
  @{ month = "May 1990"; v1 = 77028.3668376652; v2 = 59.2509705098014;  ht = 16.5 },
  @{ month = "Jun 1990"; v1 = 77298.8432535479; v2 = 59.57198969987;    ht = 16.5 },
  @{ month = "Jul 1990"; v1 = 77906.2866468195; v2 = 59.8012891213475;  ht = 16.5 },
  @{ month = "Aug 1990"; v1 = 78908.8383573794; v2 = 60.3516077328937;  ht = 16.5 },
  @{ month = "Sep 1990"; v1 = 80024.1661328505; v2 = 60.8560664601443;  ht = 16.5 },
  @{ month = "Oct 1990"; v1 = 82852.1965317919; v2 = 61.2229455345084;  ht = 16.5 },
  @{ month = "Nov 1990"; v1 = 83059.4027736163; v2 = 61.3605251873949;  ht = 16.5 },
  @{ month = "Dec 1990"; v1 = 83316.2146078445; v2 = 61.3605251873949;  ht = 17.25 }
)

$firstNewRow = 364
$lastExistingRow = 363
$lastExistingOddRow = 362   # "Nov 1989" -> banded (gray) style, alternates with 363

for ($i = 0; $i -lt $monthsData.Count; $i++) {
  $row = $firstNewRow + $i
  $data = $monthsData[$i]

  # Clone the formatting of the template row two rows back so the alternating
  # row-banding style (plain / gray) continues seamlessly into the new data.
  $templateRow = $lastExistingOddRow + ($i % 2)
  $ws1.Range("A" + $templateRow + ":D" + $templateRow).Copy()
  $ws1.Range("A" + $row + ":D" + $row).PasteSpecial(-4122)

  $ws1.Cells.Item($row, 1).Value = "United States"

  # Typing "Mon YYYY" into a General cell gets auto-parsed as a date, so force
  # the cell to Text first, enter the value, then restore its real style
  # (copied from the equivalent month cell in the template row) so no stray
  # date format lingers on the cell.
  $ws1.Cells.Item($row, 2).NumberFormat = "@"
  $ws1.Cells.Item($row, 2).Value = $data.month
  $ws1.Cells.Item($templateRow, 2).Copy()
  $ws1.Cells.Item($row, 2).PasteSpecial(-4122)

  $ws1.Cells.Item($row, 3).Value = $data.v1
  $ws1.Cells.Item($row, 4).Value = $data.v2

  $ws1.Rows.Item($row).RowHeight = $data.ht
}

$lastNewRow = $firstNewRow + $monthsData.Count - 1   # 375

# 3) The merged "United States" label in column A now spans through the new
#    last data row.
$ws1.Range("A3:A" + $lastExistingRow).UnMerge()
$ws1.Range("A3:A" + $lastNewRow).Merge()

# ---------------------------------------------------------------------------
# Sheet 2: "Tooltip"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$tooltipFirstNewRow = 363
$tooltipTemplateRow = 362   # "Dec 1989" tooltip row, used as a formatting/value template

for ($i = 0; $i -lt $monthsData.Count; $i++) {
  $row = $tooltipFirstNewRow + $i
  $data = $monthsData[$i]
  $oldMonth = "Dec 1989"
  if ($i -gt 0) { $oldMonth = $monthsData[$i - 1].month }
  $prevRow = $tooltipTemplateRow + $i

  # Clone the previous tooltip row wholesale (values + formats): country name,
  # month text and the two "Country: .. / Time: .." tooltip strings.
  $ws2.Range("A" + $prevRow + ":D" + $prevRow).Copy()
  $ws2.Range("A" + $row + ":D" + $row).PasteSpecial(-4104)

  # Fix up the free-text tooltip columns in place -- plain text substitution,
  # so it is never re-interpreted as anything other than a string.
  $ws2.Range("C" + $row + ":D" + $row).Replace($oldMonth, $data.month)

  # The month column again needs the Text-format workaround to avoid being
  # auto-converted into a date serial, then the style is reset back to the
  # sheet's default (style 0, same as every other Tooltip cell).
  $ws2.Cells.Item($row, 2).NumberFormat = "@"
  $ws2.Cells.Item($row, 2).Value = $data.month
  $ws2.Range("A1").Copy()
  $ws2.Cells.Item($row, 2).PasteSpecial(-4122)

  # Tooltip rows don't carry an explicit custom height; undo the auto-fit
  # bump that editing the wrapped tooltip text triggers.
  $ws2.Rows.Item($row).AutoFit()
}

Write-Output "edit complete"
